$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 84.76851766666665
$ws.Range("H2").Value = 254.305553
$ws.Range("I2").Value = 0.2571740874301185
$ws.Range("J2").Value = 0.2571740874301185
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 69.03718366666666
$ws.Range("N2").Value = 207.111551
$ws.Range("O2").Value = 0.7412254785079075
$ws.Range("P2").Value = 0.7412254785079075
$ws.Range("Q2").Value = 5852.179723304744
$ws.Range("R2").Value = 52669.6175097427
$ws.Range("S2").Value = 0.190623986015224
$ws.Range("T2").Value = 0.1906239860152241
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 84.76851766666665
$ws.Range("H3").Value = 254.305553
$ws.Range("I3").Value = 0.2571740874301185
$ws.Range("J3").Value = 0.2571740874301185
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.153561
$ws.Range("N3").Value = 39.460683
$ws.Range("O3").Value = 0.1412246854301422
$ws.Range("P3").Value = 0.1412246854301422
$ws.Range("Q3").Value = 1115.007868008078
$ws.Range("R3").Value = 10035.0708120727
$ws.Range("S3").Value = 0.03631932959810238
$ws.Range("T3").Value = 0.03631932959810238
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 84.76851766666665
$ws.Range("H4").Value = 254.305553
$ws.Range("I4").Value = 0.2571740874301185
$ws.Range("J4").Value = 0.2571740874301185
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.94850333333333
$ws.Range("N4").Value = 32.84551
$ws.Range("O4").Value = 0.1175498360619503
$ws.Range("P4").Value = 0.1175498360619503
$ws.Range("Q4").Value = 928.0883982352254
$ws.Range("R4").Value = 8352.795584117028
$ws.Range("S4").Value = 0.03023077181679211
$ws.Range("T4").Value = 0.03023077181679211
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 133.6830266666667
$ws.Range("H5").Value = 401.04908
$ws.Range("I5").Value = 0.4055728628296552
$ws.Range("J5").Value = 0.4055728628296552
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 69.03718366666666
$ws.Range("N5").Value = 207.111551
$ws.Range("O5").Value = 0.7412254785079075
$ws.Range("P5").Value = 0.7412254785079075
$ws.Range("Q5").Value = 9229.099665102565
$ws.Range("R5").Value = 83061.89698592307
$ws.Range("S5").Value = 0.3006209393207331
$ws.Range("T5").Value = 0.3006209393207331
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 133.6830266666667
$ws.Range("H6").Value = 401.04908
$ws.Range("I6").Value = 0.4055728628296552
$ws.Range("J6").Value = 0.4055728628296552
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.153561
$ws.Range("N6").Value = 39.460683
$ws.Range("O6").Value = 0.1412246854301422
$ws.Range("P6").Value = 0.1412246854301422
$ws.Range("Q6").Value = 1758.407845924627
$ws.Range("R6").Value = 15825.67061332164
$ws.Range("S6").Value = 0.05727689997212027
$ws.Range("T6").Value = 0.05727689997212028
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 133.6830266666667
$ws.Range("H7").Value = 401.04908
$ws.Range("I7").Value = 0.4055728628296552
$ws.Range("J7").Value = 0.4055728628296552
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.94850333333333
$ws.Range("N7").Value = 32.84551
$ws.Range("O7").Value = 0.1175498360619503
$ws.Range("P7").Value = 0.1175498360619503
$ws.Range("Q7").Value = 1463.629063070089
$ws.Range("R7").Value = 13172.6615676308
$ws.Range("S7").Value = 0.04767502353680182
$ws.Range("T7").Value = 0.04767502353680183
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 111.16377
$ws.Range("H8").Value = 333.49131
$ws.Range("I8").Value = 0.3372530497402263
$ws.Range("J8").Value = 0.3372530497402264
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 69.03718366666666
$ws.Range("N8").Value = 207.111551
$ws.Range("O8").Value = 0.7412254785079075
$ws.Range("P8").Value = 0.7412254785079075
$ws.Range("Q8").Value = 7674.43360656909
$ws.Range("R8").Value = 69069.9024591218
$ws.Range("S8").Value = 0.2499805531719504
$ws.Range("T8").Value = 0.2499805531719504
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 111.16377
$ws.Range("H9").Value = 333.49131
$ws.Range("I9").Value = 0.3372530497402263
$ws.Range("J9").Value = 0.3372530497402264
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.153561
$ws.Range("N9").Value = 39.460683
$ws.Range("O9").Value = 0.1412246854301422
$ws.Range("P9").Value = 0.1412246854301422
$ws.Range("Q9").Value = 1462.19942968497
$ws.Range("R9").Value = 13159.79486716473
$ws.Range("S9").Value = 0.04762845585991957
$ws.Range("T9").Value = 0.04762845585991957
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 111.16377
$ws.Range("H10").Value = 333.49131
$ws.Range("I10").Value = 0.3372530497402263
$ws.Range("J10").Value = 0.3372530497402264
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.94850333333333
$ws.Range("N10").Value = 32.84551
$ws.Range("O10").Value = 0.1175498360619503
$ws.Range("P10").Value = 0.1175498360619503
$ws.Range("Q10").Value = 1217.0769063909
$ws.Range("R10").Value = 10953.6921575181
$ws.Range("S10").Value = 0.03964404070835637
$ws.Range("T10").Value = 0.03964404070835638
